$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.182.96'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '2.173.23'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'237.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.45%  '
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("E7").Value = '  -4.71%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = "'0.572"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.93%  '
$ws.Range("D10").Value = "'39.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.36%  '
$ws.Range("D11").Value = "'0.0921"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.56%  '
$ws.Range("D12").Value = "'54.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.22%  '
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("E14").Value = '  -5.11%  '
$ws.Range("D15").Value = '2.501.66'
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("D16").Value = "'14.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '2.177.14'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = "'0.790"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.66%  '
$ws.Range("D19").Value = '41.080.87'
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("D20").Value = "'0.0000100"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.95%  '
$ws.Range("D21").Value = "'70.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.12%  '
$ws.Range("D22").Value = "'5.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.05%  '
$ws.Range("D23").Value = "'226.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("D24").Value = "'9.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.30%  '
$ws.Range("D25").Value = "'1.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.69%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = "'10.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.01%  '
$ws.Range("E28").Value = '  -3.87%  '
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("D31").Value = "'168.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("D32").Value = "'19.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").Value = "'30.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.30%  '
$ws.Range("D34").Value = "'0.0758"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.64%  '
$ws.Range("E35").Value = '  -9.42%  '
$ws.Range("E36").Value = '  -3.53%  '
$ws.Range("D37").Value = "'0.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.42%  '
$ws.Range("E38").Value = '  -4.08%  '
$ws.Range("D39").Value = "'0.0281"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.39%  '
$ws.Range("E40").Value = '  -2.35%  '
$ws.Range("D41").Value = "'11.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -12.44%  '
$ws.Range("D42").Value = "'5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.27%  '
$ws.Range("D43").Value = "'58.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -11.39%  '
$ws.Range("E44").Value = '  -4.34%  '
$ws.Range("E45").Value = '  -4.85%  '
$ws.Range("D46").Value = "'0.0958"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.41%  '
$ws.Range("D47").Value = "'97.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.33%  '
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("E49").Value = '  -3.54%  '
$ws.Range("E50").Value = '  -8.31%  '
$ws.Range("E51").Value = '  -2.76%  '
